$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs / Insl3 / Rxfp2 / ECs  (values recomputed with new TPM)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Insl3"
$ws.Range("C2").Value = "Rxfp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8375279999999999
$ws.Range("H2").Value = 2.512584
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.0002143333333333333
$ws.Range("N2").Value = 0.000643
$ws.Range("O2").Value = 0.0005591790590486129
$ws.Range("P2").Value = 0.0005591790590486129
$ws.Range("Q2").Value = 0.000179510168
$ws.Range("R2").Value = 0.001615591512
$ws.Range("S2").Value = 0.0005591790590486129
$ws.Range("T2").Value = 0.0005591790590486129

# Row 3 (new row): ECs / Insl3 / Rxfp2 / FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Insl3"
$ws.Range("C3").Value = "Rxfp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8375279999999999
$ws.Range("H3").Value = 2.512584
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.380202
$ws.Range("N3").Value = 1.140606
$ws.Range("O3").Value = 0.991917558048526
$ws.Range("P3").Value = 0.991917558048526
$ws.Range("Q3").Value = 0.318429820656
$ws.Range("R3").Value = 2.865868385904
$ws.Range("S3").Value = 0.991917558048526
$ws.Range("T3").Value = 0.991917558048526

# Row 4 (shifted former row 3): ECs / Insl3 / Rxfp2 / MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Insl3"
$ws.Range("C4").Value = "Rxfp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8375279999999999
$ws.Range("H4").Value = 2.512584
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002883666666666667
$ws.Range("N4").Value = 0.008651000000000001
$ws.Range("O4").Value = 0.007523262892425429
$ws.Range("P4").Value = 0.00752326289242543
$ws.Range("Q4").Value = 0.002415151576
$ws.Range("R4").Value = 0.021736364184
$ws.Range("S4").Value = 0.007523262892425429
$ws.Range("T4").Value = 0.00752326289242543
